$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6
$ws.Range("D6").Value = "[Object Detection] YOLO Define Optimal Anchor Box :: YOLO v5, YOLO v6 autoanchor"

# Row 36
$ws.Range("D36").Value = "How to Transfer Knowledge Across Domains by Deep Neural Network?"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/384"

# Row 39
$ws.Range("D39").Value = "deadNstreet"

# Row 42
$ws.Range("D42").Value = "IT_notepad"

# Row 43
$ws.Range("D43").Value = "동신한의 조재성"

# Row 44
$ws.Range("D44").Value = "Object Detection Algorithm (Efficientdet)"

# Row 46
$ws.Range("D46").Value = "[유한양행] 2022년 10월, 생물정보학(Bioinformatics 채용), 중앙연구소 인공지능 활용 신약개발 연구원"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/495"
